# Fill in the Deposition and Depletion Parameters (dep/depl/phase/pdep/pdepl/vdep/vdepl)
# for the data row (row 3): dep="Y", depl="Y", phase="B", pdep/pdepl/vdep/vdepl="WD"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3:R3").Value = "Y"
$ws.Range("S3").Value = "B"
$ws.Range("T3:W3").Value = "WD"

# Move the active selection to X3, matching the author's latest click in the sheet
$ws.Range("X3").Select()
